# Saving one hot encoder and dataset mean
# Insert two new rows above the "Scaler X file path" row on the
# HandyML_Predictor sheet, carrying the same header formatting, and
# populate them with the two new field labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HandyML_Predictor")

# Insert two blank rows before the current row 3 ("Scaler X file path"),
# shifting the existing rows (and the trailing spacer row) down.
$ws.Rows("3:4").Insert()

# Copy the formatting (style/border/font) from the still-intact
# "Scaler X file path" row (now row 5) onto the two freshly inserted rows.
$ws.Range("A5").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$ws.Range("A3:A4").RowHeight = 21

# Fill in the new labels. "One hot encoder path" is entered first so it
# lands at shared-string index 6 and "Dataset mean path" at index 7,
# matching the order the strings were appended to sharedStrings.xml.
$ws.Range("A4").Value = "One hot encoder path"
$ws.Range("A3").Value = "Dataset mean path"

# Match the author's final selection.
[void]$ws.Range("B8").Select()
